$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: find the Address of the existing hyperlink anchored at a given A1
# cell reference on a worksheet (so new hyperlinks can reuse the same target
# URL as their "sibling" handoff-file / handback-file links).
# ---------------------------------------------------------------------------
function Get-HyperlinkAddress($ws, [string]$cellRef) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq "`$$($cellRef.Substring(0,1))`$$($cellRef.Substring(1))") {
            return $hl.Address
        }
    }
    return $null
}

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: refresh the rolled-up Status column for both language rows.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# zh-cn sheet: mark handback complete, stamp the handback datetime, and fill
# in the new "Latest Target File" / "Latest Handback File" columns (F, G)
# with hyperlinked file names - mirroring the existing handoff (A) / xlf (D)
# hyperlinks.
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zhMdAddr = Get-HyperlinkAddress $zh "A2"
$zhXlfAddr = Get-HyperlinkAddress $zh "D2"

$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

$zh.Hyperlinks.Add($zh.Range("F2"), $zhMdAddr, "", "", "490b1726-4faf-4619-b87d-251fc04d19b3.md")
$zh.Hyperlinks.Add($zh.Range("G2"), $zhXlfAddr, "", "", "490b1726-4faf-4619-b87d-251fc04d19b3.5c58728acdedc51e4ff23cc4322e5fae940ad35a.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("F3"), $zhMdAddr, "", "", "490b1726-4faf-4619-b87d-251fc04d19b3.md")
$zh.Hyperlinks.Add($zh.Range("G3"), $zhXlfAddr, "", "", "490b1726-4faf-4619-b87d-251fc04d19b3.5c58728acdedc51e4ff23cc4322e5fae940ad35a.zh-cn.xlf")

$zh.Range("F2").Font.Underline = 2
$zh.Range("F2").Font.Color = 15570276
$zh.Range("G2").Font.Underline = 2
$zh.Range("G2").Font.Color = 15570276
$zh.Range("F3").Font.Underline = 2
$zh.Range("F3").Font.Color = 15570276
$zh.Range("G3").Font.Underline = 2
$zh.Range("G3").Font.Color = 15570276

$zh.Range("H2").Value = "2016-03-20 10:50:15"
$zh.Range("H3").Value = "2016-03-20 10:50:15"

# ---------------------------------------------------------------------------
# de-de sheet: same treatment, different handback timestamp / xlf file.
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$deMdAddr = Get-HyperlinkAddress $de "A2"
$deXlfAddr = Get-HyperlinkAddress $de "D2"

$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

$de.Hyperlinks.Add($de.Range("F2"), $deMdAddr, "", "", "490b1726-4faf-4619-b87d-251fc04d19b3.md")
$de.Hyperlinks.Add($de.Range("G2"), $deXlfAddr, "", "", "490b1726-4faf-4619-b87d-251fc04d19b3.5c58728acdedc51e4ff23cc4322e5fae940ad35a.de-de.xlf")
$de.Hyperlinks.Add($de.Range("F3"), $deMdAddr, "", "", "490b1726-4faf-4619-b87d-251fc04d19b3.md")
$de.Hyperlinks.Add($de.Range("G3"), $deXlfAddr, "", "", "490b1726-4faf-4619-b87d-251fc04d19b3.5c58728acdedc51e4ff23cc4322e5fae940ad35a.de-de.xlf")

$de.Range("F2").Font.Underline = 2
$de.Range("F2").Font.Color = 15570276
$de.Range("G2").Font.Underline = 2
$de.Range("G2").Font.Color = 15570276
$de.Range("F3").Font.Underline = 2
$de.Range("F3").Font.Color = 15570276
$de.Range("G3").Font.Underline = 2
$de.Range("G3").Font.Color = 15570276

$de.Range("H2").Value = "2016-03-20 10:50:20"
$de.Range("H3").Value = "2016-03-20 10:50:20"
